# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder: Kuwait now ranks above Belgica in the shared-strings table.
#     Row 38 becomes Kuwait (fresh data) and row 39 becomes Belgica
#     (keeps its previous data, just moved down one row). ---
$ws.Range("A38").Value = "Kuwait"
$ws.Range("A39").Value = "Belgica"

# --- Reorder: Islas Malvinas now listed before Montserrat. Row 214 becomes
#     Islas Malvinas and row 215 becomes Montserrat (each keeps its own
#     original data, which swaps rows along with the name). ---
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("A215").Value = "Montserrat"

# --- Refreshed case numbers (columns B:H = Casos totales, Nuevos casos,
#     Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 6877617
$ws.Range("C4").Value = 3021
$ws.Range("D4").Value = 4155933
$ws.Range("E4").Value = 2519418
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 53
$ws.Range("H4").Value = 202266

# Row 5 - India
$ws.Range("B5").Value = 5219211
$ws.Range("C5").Value = 6525
$ws.Range("D5").Value = 4112551
$ws.Range("E5").Value = 1022226
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 84434

# Row 19 - Arabia Saudita
$ws.Range("B19").Value = 328720
$ws.Range("C19").Value = 576
$ws.Range("D19").Value = 308352
$ws.Range("E19").Value = 15938
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 31
$ws.Range("H19").Value = 4430

# Row 25 - Alemania
$ws.Range("B25").Value = 269995
$ws.Range("C25").Value = 953
$ws.Range("D25").Value = 241300
$ws.Range("E25").Value = 19233
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 9462

# Row 33 - Rumania
$ws.Range("B33").Value = 110217
$ws.Range("C33").Value = 1527
$ws.Range("D33").Value = 43244
$ws.Range("E33").Value = 62613
$ws.Range("F33").Value = 0
$ws.Range("G33").Value = 48
$ws.Range("H33").Value = 4360

# Row 38 - Kuwait (new rank, fresh data)
$ws.Range("B38").Value = 98528
$ws.Range("C38").Value = 704
$ws.Range("D38").Value = 88776
$ws.Range("E38").Value = 9172
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 5
$ws.Range("H38").Value = 580

# Row 39 - Belgica (unchanged data, shifted down)
$ws.Range("B39").Value = 97976
$ws.Range("C39").Value = 2028
$ws.Range("D39").Value = 18854
$ws.Range("E39").Value = 69186
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 1
$ws.Range("H39").Value = 9936

# Row 42 - Paises Bajos
$ws.Range("B42").Value = 90047
$ws.Range("C42").Value = 1974
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 7
$ws.Range("H42").Value = 6273

# Row 43 - Suecia
$ws.Range("B43").Value = 88237
$ws.Range("C43").Value = 0
$ws.Range("D43").Value = 0
$ws.Range("E43").Value = 0
$ws.Range("F43").Value = 0
$ws.Range("G43").Value = 5
$ws.Range("H43").Value = 5865

# Row 60 - Suiza
$ws.Range("B60").Value = 49283
$ws.Range("C60").Value = 488
$ws.Range("D60").Value = 39900
$ws.Range("E60").Value = 7339
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 2
$ws.Range("H60").Value = 2044

# Row 71 - Estado de Palestina
$ws.Range("B71").Value = 34401
$ws.Range("C71").Value = 558
$ws.Range("D71").Value = 23333
$ws.Range("E71").Value = 10818
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 6
$ws.Range("H71").Value = 250

# Row 72 - Serbia
$ws.Range("B72").Value = 32757
$ws.Range("C72").Value = 62
$ws.Range("D72").Value = 31512
$ws.Range("E72").Value = 506
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 1
$ws.Range("H72").Value = 739

# Row 78 - Libia
$ws.Range("B78").Value = 26438
$ws.Range("C78").Value = 616
$ws.Range("D78").Value = 14207
$ws.Range("E78").Value = 11813
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 13
$ws.Range("H78").Value = 418

# Row 85 - Republica de Macedonia
$ws.Range("B85").Value = 16417
$ws.Range("C85").Value = 143
$ws.Range("D85").Value = 13732
$ws.Range("E85").Value = 2002
$ws.Range("F85").Value = 0
$ws.Range("G85").Value = 8
$ws.Range("H85").Value = 683

# Row 87 - Madagascar
$ws.Range("B87").Value = 15971
$ws.Range("C87").Value = 46
$ws.Range("D87").Value = 14587
$ws.Range("E87").Value = 1167
$ws.Range("F87").Value = 0
$ws.Range("G87").Value = 1
$ws.Range("H87").Value = 217

# Row 101 - Tayikistan
$ws.Range("B101").Value = 9259
$ws.Range("C101").Value = 45
$ws.Range("D101").Value = 8026
$ws.Range("E101").Value = 1160

# Row 168 - Vietnam
$ws.Range("B168").Value = 1068
$ws.Range("C168").Value = 2
$ws.Range("D168").Value = 941
$ws.Range("E168").Value = 92

# Row 214 - Islas Malvinas (data from old row 215)
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Row 215 - Montserrat (data from old row 214)
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1


# --- Update "last refreshed" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 15:21"
